# Reverse the order of the worksheet tabs.
# Target order (left to right): 总计, 2022-Q2, 2022-Q1, 2021-Q3, 2021-Q2, 2021-Q1, 2020-Q4
$wb = $excel.ActiveWorkbook

$targetOrder = @("总计", "2022-Q2", "2022-Q1", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")

for ($i = 0; $i -lt $targetOrder.Length; $i++) {
    $name = $targetOrder[$i]
    $ws = $wb.Worksheets.Item($name)
    if ($i -eq 0) {
        $ws.Move($wb.Worksheets.Item(1))
    } else {
        $prev = $wb.Worksheets.Item($targetOrder[$i - 1])
        $ws.Move($null, $prev)
    }
}
